$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Completed")

# --- Fill in the missing Finish Date for "Best. Movie. Year. Ever." (row 95) ---
# Copy the existing date-formatted style from C95 so D95 matches the sheet's
# established date-cell formatting (numFmtId 14) instead of creating a new style.
$ws.Range("C95").Copy()
$ws.Range("D95").PasteSpecial(-4122)
$ws.Range("D95").Value = 44011

# --- Add new row 97: "A Wizard of Earthsea" ---
$ws.Range("A97").Value = "A Wizard of Earthsea"
$ws.Range("B97").Value = "Ursula K Le Guin"

$ws.Range("C96").Copy()
$ws.Range("C97").PasteSpecial(-4122)
$ws.Range("C97").Value = 44011

$ws.Range("D96").Copy()
$ws.Range("D97").PasteSpecial(-4122)
$ws.Range("D97").Value = 44012

$ws.Range("E97").Value = "fiction;wizard;fantasy;quest"
$ws.Range("F97").Value = "Hard Copy"
$ws.Range("G97").Value = "145 Pages"

$ws.Range("A98").Select()
